$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance (column D) for the additional session/week for the
# students that were present, matching rows for week "sapt 3" (col D).
$ws.Range("D4").Value = $true
$ws.Range("D8").Value = $true
$ws.Range("D9").Value = $true
$ws.Range("D11").Value = $true
$ws.Range("D15").Value = $true
$ws.Range("D20").Value = $true
$ws.Range("D21").Value = $true

# Restore the selection left by the author when saving the workbook.
$ws.Range("D10").Select() | Out-Null
